$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "Pinal"
$ws.Range("C2").Value = "Raja"

$ws.Range("C2").Select()
